$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.959.10"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.588.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("E9").Value = "  -1.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.809.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.585.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.02%  "

$ws.Range("E14").Value = "  -1.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.509"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.954.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0721"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "

$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "199.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.94%  "

$ws.Range("E21").Value = "  +0.71%  "

$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  -8.39%  "

$ws.Range("E28").Value = "  -0.46%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("E31").Value = "  +0.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.33%  "

$ws.Range("E34").Value = "  -1.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.122.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.05%  "

$ws.Range("E37").Value = "  +8.61%  "

$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("E39").Value = "  -1.57%  "

$ws.Range("E40").Value = "  +0.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.488"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.776"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.21%  "

$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("E44").Value = "  -1.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.10%  "

$ws.Range("E46").Value = "  -1.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.407"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₇0929"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -17.12%  "
